# RPA datasets push 2024-03-19
#
# The "이노그리드" IPO listing had its demand-forecast (수요예측일) window
# pushed out from 2024.03.12~03.18 to 2024.04.18~04.24. On the tracking
# sheet, rows are kept sorted by that date (most recent first), so the
# row needs to move from its old spot (row 10) up to row 3, with the
# date cell updated to the new schedule. Every row that used to sit
# between the old and new position shifts down by one; everything else
# on the sheet is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02_38커뮤니케이션(최근일자기준)")

# Make room for the row at its new position (row 3) - shifts rows 3-21 down to 4-22.
$ws.Rows("3:3").Insert()

# Write the 이노그리드 row into its new spot with the updated demand-forecast date.
$ws.Range("A3").Value = "이노그리드"
$ws.Range("B3").Value = "2024.04.18~04.24"
$ws.Range("C3").Value = "29,000~35,000"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 17400
$ws.Range("F3").Value = "한국투자증권"

# Remove the row's old location (now pushed down to row 11 by the insert above),
# shifting rows 12-22 back up to 11-21 so the rest of the table is unchanged.
$ws.Rows("11:11").Delete()
